# Scheduled-runner style market-data refresh: updates price/profit columns
# (H..N) on several leve-profit sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 276.66666
$ws.Range("I2").Value = 276.66666
$ws.Range("K2").Value = 276.66666
$ws.Range("M2").Value = -163.66666
$ws.Range("H9").Value = 226.66667
$ws.Range("I9").Value = 170
$ws.Range("J9").Value = 255
$ws.Range("K9").Value = 170
$ws.Range("L9").Value = 255
$ws.Range("M9").Value = -1
$ws.Range("N9").Value = -593
$ws.Range("I43").Value = 417.33334
$ws.Range("J43").Value = 1459.091
$ws.Range("K43").Value = 417.33334
$ws.Range("L43").Value = 1459.091
$ws.Range("M43").Value = -348.33334
$ws.Range("N43").Value = -1597.091
$ws.Range("H116").Value = 2623.2942
$ws.Range("J116").Value = 3199.2
$ws.Range("L116").Value = 3199.2
$ws.Range("N116").Value = -10083.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9520.8125
$ws.Range("I61").Value = 8288
$ws.Range("J61").Value = 11105.857
$ws.Range("K61").Value = 8288
$ws.Range("L61").Value = 11105.857
$ws.Range("M61").Value = -8076
$ws.Range("N61").Value = -11529.857
$ws.Range("H136").Value = 9520.8125
$ws.Range("I136").Value = 8288
$ws.Range("J136").Value = 11105.857
$ws.Range("K136").Value = 24864
$ws.Range("L136").Value = 33317.571
$ws.Range("M136").Value = -22314
$ws.Range("N136").Value = -38417.571

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 955.2
$ws.Range("I5").Value = 178
$ws.Range("J5").Value = 1473.3334
$ws.Range("K5").Value = 178
$ws.Range("L5").Value = 1473.3334
$ws.Range("M5").Value = -66
$ws.Range("N5").Value = -1697.3334
$ws.Range("H25").Value = 12400
$ws.Range("J25").Value = 12400
$ws.Range("L25").Value = 12400
$ws.Range("N25").Value = -12748
$ws.Range("H53").Value = 34710
$ws.Range("J53").Value = 34710
$ws.Range("L53").Value = 34710
$ws.Range("N53").Value = -35924
$ws.Range("H58").Value = 2167091.2
$ws.Range("I58").Value = 3498185.5
$ws.Range("J58").Value = 4063.25
$ws.Range("K58").Value = 3498185.5
$ws.Range("L58").Value = 4063.25
$ws.Range("M58").Value = -3497982.5
$ws.Range("N58").Value = -4469.25
$ws.Range("H111").Value = 74900
$ws.Range("J111").Value = 74900
$ws.Range("L111").Value = 74900
$ws.Range("N111").Value = -83080
$ws.Range("H132").Value = 3310.9565
$ws.Range("I132").Value = 3161.5557
$ws.Range("J132").Value = 3407
$ws.Range("K132").Value = 9484.667099999999
$ws.Range("L132").Value = 10221
$ws.Range("M132").Value = -6954.667099999999
$ws.Range("N132").Value = -15281
$ws.Range("H134").Value = 2823.2654
$ws.Range("I134").Value = 2226
$ws.Range("K134").Value = 6678
$ws.Range("M134").Value = -4143
$ws.Range("H136").Value = 2167091.2
$ws.Range("I136").Value = 3498185.5
$ws.Range("J136").Value = 4063.25
$ws.Range("K136").Value = 10494556.5
$ws.Range("L136").Value = 12189.75
$ws.Range("M136").Value = -10492006.5
$ws.Range("N136").Value = -17289.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 849.75
$ws.Range("I7").Value = 833
$ws.Range("K7").Value = 2499
$ws.Range("M7").Value = -2387
$ws.Range("H8").Value = 72.333336
$ws.Range("I8").Value = 72.333336
$ws.Range("K8").Value = 217.000008
$ws.Range("M8").Value = -78.00000800000001
$ws.Range("H22").Value = 2533.1667
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 1800
$ws.Range("M22").Value = -1631
$ws.Range("H27").Value = 2533.1667
$ws.Range("I27").Value = 600
$ws.Range("K27").Value = 1800
$ws.Range("M27").Value = -1698
$ws.Range("H33").Value = 212.85715
$ws.Range("I33").Value = 246.27272
$ws.Range("J33").Value = 90.333336
$ws.Range("K33").Value = 1477.63632
$ws.Range("L33").Value = 542.000016
$ws.Range("M33").Value = -1194.63632
$ws.Range("N33").Value = -1108.000016
$ws.Range("H74").Value = 5999
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = 17997
$ws.Range("N74").Value = -20119
$ws.Range("L74").ClearContents()
$ws.Range("H77").Value = 5999
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = 53991
$ws.Range("N77").Value = -64599
$ws.Range("L77").ClearContents()
$ws.Range("H80").Value = 2400.182
$ws.Range("I80").Value = 2025.5
$ws.Range("J80").Value = 2614.2856
$ws.Range("K80").Value = 6076.5
$ws.Range("L80").Value = 7842.8568
$ws.Range("M80").Value = -5140.5
$ws.Range("N80").Value = -9714.856800000001
$ws.Range("H83").Value = 2400.182
$ws.Range("I83").Value = 2025.5
$ws.Range("J83").Value = 2614.2856
$ws.Range("K83").Value = 18229.5
$ws.Range("L83").Value = 23528.5704
$ws.Range("M83").Value = -13549.5
$ws.Range("N83").Value = -32888.5704
$ws.Range("H86").Value = 620.55554
$ws.Range("I86").Value = 641
$ws.Range("J86").Value = 595
$ws.Range("K86").Value = 1923
$ws.Range("L86").Value = 1785
$ws.Range("M86").Value = -737
$ws.Range("N86").Value = -4157
$ws.Range("H89").Value = 620.55554
$ws.Range("I89").Value = 641
$ws.Range("J89").Value = 595
$ws.Range("K89").Value = 5769
$ws.Range("L89").Value = 5355
$ws.Range("M89").Value = 159
$ws.Range("N89").Value = -17211
$ws.Range("H92").Value = 1196.5
$ws.Range("J92").Value = 1403
$ws.Range("L92").Value = 4209
$ws.Range("N92").Value = -6705
$ws.Range("H97").Value = 13261.111
$ws.Range("I97").Value = 1287.5
$ws.Range("K97").Value = 3862.5
$ws.Range("M97").Value = -3366.5
$ws.Range("H107").Value = 1055.8667
$ws.Range("J107").Value = 1508.8889
$ws.Range("L107").Value = 4526.6667
$ws.Range("N107").Value = -8366.6667
$ws.Range("H116").Value = 500
$ws.Range("I116").Value = 500
$ws.Range("K116").Value = 1500
$ws.Range("M116").Value = 1942
$ws.Range("H122").Value = 965.0454999999999
$ws.Range("J122").Value = 1007.8333
$ws.Range("L122").Value = 9070.4997
$ws.Range("N122").Value = -13970.4997
$ws.Range("H125").Value = 2691
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 2799.2727
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 8397.8181
$ws.Range("M125").Value = 420
$ws.Range("N125").Value = -18237.8181
$ws.Range("H132").Value = 2241.7856
$ws.Range("J132").Value = 1691.3636
$ws.Range("L132").Value = 15222.2724
$ws.Range("N132").Value = -20282.2724

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 27500
$ws.Range("J94").Value = 27500
$ws.Range("L94").Value = 27500
$ws.Range("N94").Value = -28852
$ws.Range("H136").Value = 5283.8
$ws.Range("I136").Value = 4347.8
$ws.Range("J136").Value = 6219.8
$ws.Range("K136").Value = 13043.4
$ws.Range("L136").Value = 18659.4
$ws.Range("M136").Value = -10493.4
$ws.Range("N136").Value = -23759.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 21608
$ws.Range("J104").Value = 21608
$ws.Range("L104").Value = 21608
$ws.Range("N104").Value = -28596
$ws.Range("H132").Value = 1732.6052
$ws.Range("I132").Value = 1032.9524
$ws.Range("K132").Value = 3098.857199999999
$ws.Range("M132").Value = -568.8571999999995
$ws.Range("H136").Value = 3253.0889
$ws.Range("I136").Value = 2504.8276
$ws.Range("K136").Value = 7514.4828
$ws.Range("M136").Value = -4964.4828
